# Add a new "Disable Flag *" column (G) to the SSE sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell G1 -------------------------------------------------
# Start from the same look as the other header cells (bold font,
# full thin border, centered/top aligned) by copying F1's formatting,
# then trim the top/bottom border so only the left/right edges remain.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$head = $ws.Range("G1")
$head.Value = "Disable Flag *"
$head.Borders.Item(8).LineStyle = -4142          # xlEdgeTop    -> none
$head.Borders.Item(9).LineStyle = -4142          # xlEdgeBottom -> none

# --- Data rows 2-198 -------------------------------------------------
# New "Disable Flag *" values: the first few records are flagged "N"
# and the remainder of the sheet is flagged "Y".
for ($r = 2; $r -le 198; $r++) {
    if ($r -le 6) {
        $ws.Cells.Item($r, 7).Value = "N"
    } else {
        $ws.Cells.Item($r, 7).Value = "Y"
    }
}

# Restore the cursor/selection to where the author last left it.
$ws.Range("E7").Select() | Out-Null
